# Refactor code structure and remove redundant sections for improved readability
#
# - Drop the per-row picture column ("标签"/tags column E) and replace the
#   DISPIMG() image formulas in column D with plain text file names,
#   renaming the column D header from "图片" to "图片名称".
# - Clear out the (now unused) cellImages part.
# - Re-point the active selection and shrink the application window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "image" column header and replace each DISPIMG(...) formula
# cell with the literal image file name it used to display.
$ws.Range("D1").Value = "图片名称"
$ws.Range("D2").Value = "xiaojinyu.png"
$ws.Range("D3").Value = "jiatelin.png"
$ws.Range("D4").Value = "paozhang.png"
$ws.Range("D5").Value = "hudie.png"

# The old "标签" (tags) column is no longer used - remove it entirely.
$ws.Columns("E:E").Delete()

# Deleting the last grouped column drops the sheet's column outline depth
# by one level (row grouping depth is untouched).
$ws.Outline.ShowLevels(4, 3)

# The DISPIMG() formulas previously forced tall rows to show the pictures;
# now that the cells just hold plain file-name text, let the row heights
# shrink back down to the sheet default.
$ws.Rows("2:5").AutoFit()

# Match the author's final selection / window sizing.
$ws.Range("E7").Select()

$win = $excel.Windows.Item(1)
$win.Width = 14400
